$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "item_reference"
$ws.Range("B1").Value = "order"

# Data rows
$ws.Range("A2").Value = 1000012164
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = 1000007500
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 1000000676
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = 1000000633
$ws.Range("B5").Value = 4

$ws.Range("A6").Value = 1000014725
$ws.Range("B6").Value = 5

# Update selection to match target
$ws.Range("A1:B6").Select()
